$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("servicos")

# Update service names (column B)
$ws.Range("B2").Value = "Corte clássico"
$ws.Range("B3").Value = "Barba e bigode"
$ws.Range("B4").Value = "Combo clássico"
$ws.Range("B5").Value = "Tratamento capilar"
$ws.Range("B6").Value = "Dia do noivo"

# Update descriptions (column C)
$ws.Range("C2").Value = "Aparagem e modelagem tradicionais"
$ws.Range("C3").Value = "Desenho e contorno da barba com toalha quente e finalização com balm"
$ws.Range("C4").Value = "Corte de cabelo e barba clássica"
$ws.Range("C5").Value = "Lavagem especial com massagem e hidratação profunda para revitalização dos fios"
$ws.Range("C6").Value = "Pacote completo com corte, barba, tratamento facial e relaxamento"

# Update codigo for row 6
$ws.Range("A6").Value = 5

# Update values (column D)
$ws.Range("D2").Value = 40
$ws.Range("D3").Value = 35
$ws.Range("D4").Value = 70
$ws.Range("D5").Value = 50
$ws.Range("D6").Value = 150
